$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 16.820675
$ws.Cells.Item(2, 8).Value = 50.462025
$ws.Cells.Item(2, 9).Value = 0.8427583848046372
$ws.Cells.Item(2, 10).Value = 0.8427583848046373
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.3930163333333334
$ws.Cells.Item(2, 14).Value = 1.179049
$ws.Cells.Item(2, 15).Value = 0.03108558724574714
$ws.Cells.Item(2, 16).Value = 0.03108558724574714
$ws.Cells.Item(2, 17).Value = 6.610800012691667
$ws.Cells.Item(2, 18).Value = 59.497200114225
$ws.Cells.Item(2, 19).Value = 0.02619763929792949
$ws.Cells.Item(2, 20).Value = 0.0261976392979295

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 16.820675
$ws.Cells.Item(3, 8).Value = 50.462025
$ws.Cells.Item(3, 9).Value = 0.8427583848046372
$ws.Cells.Item(3, 10).Value = 0.8427583848046373
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 6.505607333333333
$ws.Cells.Item(3, 14).Value = 19.516822
$ws.Cells.Item(3, 15).Value = 0.5145603558806437
$ws.Cells.Item(3, 16).Value = 0.5145603558806437
$ws.Cells.Item(3, 17).Value = 109.4287066316166
$ws.Cells.Item(3, 18).Value = 984.8583596845498
$ws.Cells.Item(3, 19).Value = 0.4336500544064706
$ws.Cells.Item(3, 20).Value = 0.4336500544064706

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 16.820675
$ws.Cells.Item(4, 8).Value = 50.462025
$ws.Cells.Item(4, 9).Value = 0.8427583848046372
$ws.Cells.Item(4, 10).Value = 0.8427583848046373
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.737319
$ws.Cells.Item(4, 14).Value = 2.211957
$ws.Cells.Item(4, 15).Value = 0.05831817193970827
$ws.Cells.Item(4, 16).Value = 0.05831817193970829
$ws.Cells.Item(4, 17).Value = 12.402203270325
$ws.Cells.Item(4, 18).Value = 111.619829432925
$ws.Cells.Item(4, 19).Value = 0.04914812838866766
$ws.Cells.Item(4, 20).Value = 0.04914812838866767

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 16.820675
$ws.Cells.Item(5, 8).Value = 50.462025
$ws.Cells.Item(5, 9).Value = 0.8427583848046372
$ws.Cells.Item(5, 10).Value = 0.8427583848046373
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 5.007097666666667
$ws.Cells.Item(5, 14).Value = 15.021293
$ws.Cells.Item(5, 15).Value = 0.3960358849339008
$ws.Cells.Item(5, 16).Value = 0.3960358849339009
$ws.Cells.Item(5, 17).Value = 84.22276254425832
$ws.Cells.Item(5, 18).Value = 758.004862898325
$ws.Cells.Item(5, 19).Value = 0.3337625627115693
$ws.Cells.Item(5, 20).Value = 0.3337625627115695

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.464483
$ws.Cells.Item(6, 8).Value = 4.393449
$ws.Cells.Item(6, 9).Value = 0.07337430439942808
$ws.Cells.Item(6, 10).Value = 0.07337430439942808
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.3930163333333334
$ws.Cells.Item(6, 14).Value = 1.179049
$ws.Cells.Item(6, 15).Value = 0.03108558724574714
$ws.Cells.Item(6, 16).Value = 0.03108558724574714
$ws.Cells.Item(6, 17).Value = 0.5755657388890001
$ws.Cells.Item(6, 18).Value = 5.180091650001001
$ws.Cells.Item(6, 19).Value = 0.00228088334100443
$ws.Cells.Item(6, 20).Value = 0.00228088334100443

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.464483
$ws.Cells.Item(7, 8).Value = 4.393449
$ws.Cells.Item(7, 9).Value = 0.07337430439942808
$ws.Cells.Item(7, 10).Value = 0.07337430439942808
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 6.505607333333333
$ws.Cells.Item(7, 14).Value = 19.516822
$ws.Cells.Item(7, 15).Value = 0.5145603558806437
$ws.Cells.Item(7, 16).Value = 0.5145603558806437
$ws.Cells.Item(7, 17).Value = 9.527351344342001
$ws.Cells.Item(7, 18).Value = 85.74616209907799
$ws.Cells.Item(7, 19).Value = 0.0377555081842644
$ws.Cells.Item(7, 20).Value = 0.0377555081842644

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.464483
$ws.Cells.Item(8, 8).Value = 4.393449
$ws.Cells.Item(8, 9).Value = 0.07337430439942808
$ws.Cells.Item(8, 10).Value = 0.07337430439942808
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.737319
$ws.Cells.Item(8, 14).Value = 2.211957
$ws.Cells.Item(8, 15).Value = 0.05831817193970827
$ws.Cells.Item(8, 16).Value = 0.05831817193970829
$ws.Cells.Item(8, 17).Value = 1.079791141077
$ws.Cells.Item(8, 18).Value = 9.718120269693001
$ws.Cells.Item(8, 19).Value = 0.00427905529992234
$ws.Cells.Item(8, 20).Value = 0.004279055299922341

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.464483
$ws.Cells.Item(9, 8).Value = 4.393449
$ws.Cells.Item(9, 9).Value = 0.07337430439942808
$ws.Cells.Item(9, 10).Value = 0.07337430439942808
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 5.007097666666667
$ws.Cells.Item(9, 14).Value = 15.021293
$ws.Cells.Item(9, 15).Value = 0.3960358849339008
$ws.Cells.Item(9, 16).Value = 0.3960358849339009
$ws.Cells.Item(9, 17).Value = 7.332809412173001
$ws.Cells.Item(9, 18).Value = 65.99528470955701
$ws.Cells.Item(9, 19).Value = 0.02905885757423691
$ws.Cells.Item(9, 20).Value = 0.02905885757423691

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.1122456666666666
$ws.Cells.Item(10, 8).Value = 0.336737
$ws.Cells.Item(10, 9).Value = 0.005623791954919746
$ws.Cells.Item(10, 10).Value = 0.005623791954919746
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.3930163333333334
$ws.Cells.Item(10, 14).Value = 1.179049
$ws.Cells.Item(10, 15).Value = 0.03108558724574714
$ws.Cells.Item(10, 16).Value = 0.03108558724574714
$ws.Cells.Item(10, 17).Value = 0.04411438034588888
$ws.Cells.Item(10, 18).Value = 0.3970294231129999
$ws.Cells.Item(10, 19).Value = 0.0001748188754665887
$ws.Cells.Item(10, 20).Value = 0.0001748188754665887

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.1122456666666666
$ws.Cells.Item(11, 8).Value = 0.336737
$ws.Cells.Item(11, 9).Value = 0.005623791954919746
$ws.Cells.Item(11, 10).Value = 0.005623791954919746
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 6.505607333333333
$ws.Cells.Item(11, 14).Value = 19.516822
$ws.Cells.Item(11, 15).Value = 0.5145603558806437
$ws.Cells.Item(11, 16).Value = 0.5145603558806437
$ws.Cells.Item(11, 17).Value = 0.7302262322015554
$ws.Cells.Item(11, 18).Value = 6.572036089813999
$ws.Cells.Item(11, 19).Value = 0.002893780389722206
$ws.Cells.Item(11, 20).Value = 0.002893780389722206

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.1122456666666666
$ws.Cells.Item(12, 8).Value = 0.336737
$ws.Cells.Item(12, 9).Value = 0.005623791954919746
$ws.Cells.Item(12, 10).Value = 0.005623791954919746
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.737319
$ws.Cells.Item(12, 14).Value = 2.211957
$ws.Cells.Item(12, 15).Value = 0.05831817193970827
$ws.Cells.Item(12, 16).Value = 0.05831817193970829
$ws.Cells.Item(12, 17).Value = 0.08276086270099998
$ws.Cells.Item(12, 18).Value = 0.7448477643089999
$ws.Cells.Item(12, 19).Value = 0.0003279692661801578
$ws.Cells.Item(12, 20).Value = 0.0003279692661801579

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.1122456666666666
$ws.Cells.Item(13, 8).Value = 0.336737
$ws.Cells.Item(13, 9).Value = 0.005623791954919746
$ws.Cells.Item(13, 10).Value = 0.005623791954919746
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 5.007097666666667
$ws.Cells.Item(13, 14).Value = 15.021293
$ws.Cells.Item(13, 15).Value = 0.3960358849339008
$ws.Cells.Item(13, 16).Value = 0.3960358849339009
$ws.Cells.Item(13, 17).Value = 0.562025015660111
$ws.Cells.Item(13, 18).Value = 5.058225140941
$ws.Cells.Item(13, 19).Value = 0.002227223423550793
$ws.Cells.Item(13, 20).Value = 0.002227223423550794

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1.561668
$ws.Cells.Item(14, 8).Value = 4.685003999999999
$ws.Cells.Item(14, 9).Value = 0.0782435188410149
$ws.Cells.Item(14, 10).Value = 0.0782435188410149
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.3930163333333334
$ws.Cells.Item(14, 14).Value = 1.179049
$ws.Cells.Item(14, 15).Value = 0.03108558724574714
$ws.Cells.Item(14, 16).Value = 0.03108558724574714
$ws.Cells.Item(14, 17).Value = 0.6137610312439999
$ws.Cells.Item(14, 18).Value = 5.523849281195999
$ws.Cells.Item(14, 19).Value = 0.002432245731346629
$ws.Cells.Item(14, 20).Value = 0.002432245731346629

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 1.561668
$ws.Cells.Item(15, 8).Value = 4.685003999999999
$ws.Cells.Item(15, 9).Value = 0.0782435188410149
$ws.Cells.Item(15, 10).Value = 0.0782435188410149
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 6.505607333333333
$ws.Cells.Item(15, 14).Value = 19.516822
$ws.Cells.Item(15, 15).Value = 0.5145603558806437
$ws.Cells.Item(15, 16).Value = 0.5145603558806437
$ws.Cells.Item(15, 17).Value = 10.159598793032
$ws.Cells.Item(15, 18).Value = 91.43638913728797
$ws.Cells.Item(15, 19).Value = 0.04026101290018648
$ws.Cells.Item(15, 20).Value = 0.04026101290018648

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 1.561668
$ws.Cells.Item(16, 8).Value = 4.685003999999999
$ws.Cells.Item(16, 9).Value = 0.0782435188410149
$ws.Cells.Item(16, 10).Value = 0.0782435188410149
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.737319
$ws.Cells.Item(16, 14).Value = 2.211957
$ws.Cells.Item(16, 15).Value = 0.05831817193970827
$ws.Cells.Item(16, 16).Value = 0.05831817193970829
$ws.Cells.Item(16, 17).Value = 1.151447488092
$ws.Cells.Item(16, 18).Value = 10.363027392828
$ws.Cells.Item(16, 19).Value = 0.004563018984938111
$ws.Cells.Item(16, 20).Value = 0.004563018984938111

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 1.561668
$ws.Cells.Item(17, 8).Value = 4.685003999999999
$ws.Cells.Item(17, 9).Value = 0.0782435188410149
$ws.Cells.Item(17, 10).Value = 0.0782435188410149
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 5.007097666666667
$ws.Cells.Item(17, 14).Value = 15.021293
$ws.Cells.Item(17, 15).Value = 0.3960358849339008
$ws.Cells.Item(17, 16).Value = 0.3960358849339009
$ws.Cells.Item(17, 17).Value = 7.819424198907999
$ws.Cells.Item(17, 18).Value = 70.37481779017199
$ws.Cells.Item(17, 19).Value = 0.03098724122454367
$ws.Cells.Item(17, 20).Value = 0.03098724122454368

